# "Generic Backlog" is the active sheet (sheet1.xml) in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 14 ("Verizon Reimbursment"),
# shifting rows 14-21 down to 15-22.
$ws.Rows("14").Insert()

# Copy the formatting (styles) of the row that just got pushed down to 15
# (originally row 14) onto the newly-inserted blank row 14, so the new
# task line matches the "Essential" section's row styling.
$ws.Range("A15:B15").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)

# Fill in the new backlog item and its status.
$ws.Range("A14").Value = "Find a forum for Personal Finance /  Questions about  Mortgage"
$ws.Range("B14").Value = "TODO"

# Match the author's final selection state.
$ws.Range("A14").Select()
